# Swap the two name-parts of "Sanjana Meena" -> "Mena Sanjana" in the
# "Team Member 1 Name" row of the team-roster textbox on slide 4, leaving
# every other run / formatting attribute untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Locate the shape that actually holds the text (the big body placeholder
# listing team leader / members / mentor).
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -like "*Sanjana Meena*") {
            $target = $sh
        }
    }
}

$tr = $target.TextFrame.TextRange
$hit = $tr.Find("Sanjana Meena", 0)
$hit.Text = "Mena Sanjana"
